# Add a newly-received/missing product row ("STARVILLE WHITENING ANTIPRESPIRANT ROLL-ON 60 ML")
# into the "Missing Items" report, keeping the existing alphabetical ordering of products,
# and refresh the report's generation timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The product list is alphabetically sorted; the new product belongs right before
# "STREPTOQUIN 20 TABLETS", which currently lives on row 55.
$newRow = 55
$afterNewRow = $newRow + 1

# Insert a blank row at position 55, shifting the existing data (and the totals/footer rows) down.
$ws.Rows("$($newRow):$($newRow)").Insert()

# Copy the formatting (styles, merged cells, row height, fonts, borders, number formats, ...)
# of the row that was just pushed down, so that the new row looks identical to its neighbours.
$srcFormat = $ws.Range("A$($afterNewRow):Q$($afterNewRow)")
$dstFormat = $ws.Range("A$($newRow):Q$($newRow)")
$srcFormat.Copy($dstFormat)
$ws.Rows($newRow).RowHeight = $ws.Rows($afterNewRow).RowHeight
$excel.CutCopyMode = 0

# Make sure the text-like columns stay stored as text (same convention used by every other
# data row), then fill in the values for the new product.
$textCols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")
foreach ($col in $textCols) {
    $ws.Range("$col$newRow").NumberFormat = "@"
}

$ws.Range("C$newRow").Value = "STARVILLE WHITENING ANTIPRESPIRANT ROLL-ON 60 ML"
$ws.Range("H$newRow").Value = "1:0"
$ws.Range("L$newRow").Value = "1"
$ws.Range("N$newRow").Value = "130.00"
$ws.Range("P$newRow").Value = "130.0000"
$ws.Range("Q$newRow").Value = "1:0"

# Renumber the "#" column sequentially for every data row, since the new product shifts the
# position of every product that sorts after it.
$firstDataRow = 7
$lastDataRow = 85
$idx = 1
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ws.Range("A$r").Value = $idx
    $idx++
}

# Update the cached grand-total (selling-price column) shown on the totals row right below the
# data: it must grow by the new product's selling price (130.0000).
$totalsRow = $lastDataRow + 1
$oldTotal = $ws.Range("P$totalsRow").Value()
$ws.Range("P$totalsRow").Value = [double]$oldTotal + 130

# Refresh the "generated at" timestamp shown in the footer row.
$footerRow = $totalsRow + 1
$ws.Range("A$footerRow").Value = "Sunday, 17 August, 2025 5:24 PM"
